# Refresh the cryptocurrency price / 1h-volume snapshot (and fix the
# Quant / NEARProtocol row ordering) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a
# number (e.g. "1.011"); force them to Text first so the literal
# string is preserved, matching the scraped inline-string data.
$textForceCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D18",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D50"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Price (D) / Volume 1h (E) updates ---
$ws.Range("D2").Value = "27.424.04"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "1.862.79"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "311.19"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "0.4774"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.3796"
$ws.Range("E8").Value = "  +3.27%  "
$ws.Range("D9").Value = "0.07320"
$ws.Range("E9").Value = "  +1.46%  "
$ws.Range("D10").Value = "0.9326"
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("D11").Value = "20.75"
$ws.Range("E11").Value = "  +5.24%  "
$ws.Range("D12").Value = "0.07807"
$ws.Range("E12").Value = "  +0.92%  "
$ws.Range("D13").Value = "1.886.43"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").Value = "5.431"
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").Value = "6.546"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "90.29"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "0.000008798"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "27.499.46"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("D21").Value = "14.65"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").Value = "5.110"
$ws.Range("E22").Value = "  +1.00%  "
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "1.942"
$ws.Range("E24").Value = "  +0.67%  "
$ws.Range("D25").Value = "155.61"
$ws.Range("E25").Value = "  +1.76%  "
$ws.Range("D26").Value = "18.48"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "115.26"
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("D29").Value = "4.955"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "0.08876"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "3.326"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("E32").Value = "  +3.51%  "
$ws.Range("D33").Value = "0.7583"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "4.600"
$ws.Range("E34").Value = "  +2.28%  "
$ws.Range("D35").Value = "2.718"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "0.02050"
$ws.Range("E36").Value = "  +4.48%  "
$ws.Range("D37").Value = "1.122"
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("D38").Value = "0.5557"
$ws.Range("E38").Value = "  +6.62%  "
$ws.Range("D39").Value = "0.05266"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "2.991"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "7.054"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").Value = "8.634"
$ws.Range("E42").Value = "  +4.75%  "
$ws.Range("E43").Value = "  +0.81%  "
$ws.Range("D44").Value = "0.4894"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").Value = "10.66"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("D47").Value = "1.659"
$ws.Range("E47").Value = "  +3.34%  "
$ws.Range("D48").Value = "102.91"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").Value = "67.44"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("D50").Value = "0.06082"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  +3.17%  "

# --- Row 47/48 content swap: NEARProtocol now ranks above Quant ---
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
